$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row above row 16, shifting the existing data
# (rows 16-37) down to rows 17-38.
$ws.Rows("16:16").Insert()

# Populate the newly inserted row 16 with the new weekly record.
$ws.Range("A16").Value = 1
$ws.Range("B16").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C16").Value = "Arica y Parinacota"
$ws.Range("D16").Value = 45264
$ws.Range("E16").Value = 15
$ws.Range("F16").Value = 100114007
$ws.Range("G16").Value = "Jengibre"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 300
$ws.Range("K16").Value = 16000
$ws.Range("L16").Value = 18000
$ws.Range("M16").Value = 17000
$ws.Range("N16").Value = "`$/caja 13 kilos"
$ws.Range("O16").Value = "Perú"
$ws.Range("P16").Value = 1308
$ws.Range("Q16").Value = 13
$ws.Range("R16").Value = "Hortaliza"
